# Adding column with versions of browsers to use in docker containers
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# DATA sheet: fill in the "version" column plus a handful of other
# corrections to the browser/account data used by the excel-driven tests.
# ---------------------------------------------------------------------
$data = $wb.Worksheets.Item("DATA")
$data.Activate()

$data.Range("B2").Value = "no"

$data.Range("C3").Value = "chrome"
$data.Range("E3").Value = "Admin1"
$data.Range("F3").Value = "admin123"

$data.Range("B4").Value = "yes"
$data.Range("D4").Value = "'95.0.2"
$data.Range("E4").Value = "Admin"
$data.Range("F4").Value = "admin123"

$data.Range("C5").Value = "firefox"
$data.Range("E5").Value = "Admin"
$data.Range("F5").Value = "admin123"

$data.Range("E6").Value = "Admin1"

$data.Range("C7").Value = "chrome"
$data.Range("E7").Value = "Admin"
$data.Range("F7").Value = "admin123"

# Selection left on D4 (single cell) after reviewing the new version column.
$data.Range("D4").Select()

# ---------------------------------------------------------------------
# RUNMANAGER sheet: no data changed, only the remembered selection.
# ---------------------------------------------------------------------
$runmanager = $wb.Worksheets.Item("RUNMANAGER")
$runmanager.Activate()
$runmanager.Range("C2:C5").Select()

# Leave DATA as the active sheet, matching the saved workbook view.
$data.Activate()
